# Insert a new price-record row at row 130 (pushing the existing rows
# 130..192 down to 131..193, matching the target diff which grows the
# sheet from A1:T192 to A1:T193).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new record's data. The
# shared columns (A,B,C,E,F,G,H,I,J,K,Q,T) carry the same constant values
# used by every other row in this subset.
$ws.Cells.Item(130, 1).Value  = 3
$ws.Cells.Item(130, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(130, 3).Value  = "Coquimbo"
$ws.Cells.Item(130, 4).Value  = 44452
$ws.Cells.Item(130, 5).Value  = 5
$ws.Cells.Item(130, 6).Value  = "Fruta"
$ws.Cells.Item(130, 7).Value  = 100108
$ws.Cells.Item(130, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(130, 9).Value  = 100108002
$ws.Cells.Item(130, 10).Value = "Mango"
$ws.Cells.Item(130, 11).Value = "Sin especificar"
$ws.Cells.Item(130, 12).Value = "Primera"
$ws.Cells.Item(130, 13).Value = 456
$ws.Cells.Item(130, 14).Value = 9000
$ws.Cells.Item(130, 15).Value = 9000
$ws.Cells.Item(130, 16).Value = 9000
$ws.Cells.Item(130, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(130, 18).Value = "Brasil"
$ws.Cells.Item(130, 19).Value = 2250
$ws.Cells.Item(130, 20).Value = 4
